$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 20:58:21"
$ws1.Cells.Item(3,1).Value = "Total filas: 662"

$ws1.Cells.Item(648,2).Value = "20:58:10"
$ws1.Cells.Item(648,3).Value = "21:04"
$ws1.Cells.Item(648,4).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(648,5).Value = 6
$ws1.Cells.Item(648,6).Value = "LP1912"
$ws1.Cells.Item(648,7).Value = "30/12/2025"

$ws1.Cells.Item(649,2).Value = "20:58:10"
$ws1.Cells.Item(649,3).Value = "21:07"
$ws1.Cells.Item(649,4).Value = "215B_EL PATO"
$ws1.Cells.Item(649,5).Value = 9
$ws1.Cells.Item(649,6).Value = "LP1912"
$ws1.Cells.Item(649,7).Value = "30/12/2025"

$ws1.Cells.Item(650,2).Value = "20:58:10"
$ws1.Cells.Item(650,3).Value = "21:18"
$ws1.Cells.Item(650,4).Value = "16_SANTA ANA"
$ws1.Cells.Item(650,5).Value = 20
$ws1.Cells.Item(650,6).Value = "LP1912"
$ws1.Cells.Item(650,7).Value = "30/12/2025"

$ws1.Cells.Item(651,2).Value = "20:58:10"
$ws1.Cells.Item(651,3).Value = "21:21"
$ws1.Cells.Item(651,4).Value = "26_HERNANDEZ"
$ws1.Cells.Item(651,5).Value = 23
$ws1.Cells.Item(651,6).Value = "LP1912"
$ws1.Cells.Item(651,7).Value = "30/12/2025"

$ws1.Cells.Item(652,2).Value = "20:58:10"
$ws1.Cells.Item(652,3).Value = "21:23"
$ws1.Cells.Item(652,4).Value = "15_ABASTO"
$ws1.Cells.Item(652,5).Value = 25
$ws1.Cells.Item(652,6).Value = "LP1912"
$ws1.Cells.Item(652,7).Value = "30/12/2025"

$ws1.Cells.Item(653,2).Value = "20:58:10"
$ws1.Cells.Item(653,3).Value = "21:25"
$ws1.Cells.Item(653,4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(653,5).Value = 27
$ws1.Cells.Item(653,6).Value = "LP1912"
$ws1.Cells.Item(653,7).Value = "30/12/2025"

$ws1.Cells.Item(654,2).Value = "20:58:10"
$ws1.Cells.Item(654,3).Value = "21:32"
$ws1.Cells.Item(654,4).Value = "16_SANTA ANA"
$ws1.Cells.Item(654,5).Value = 34
$ws1.Cells.Item(654,6).Value = "LP1912"
$ws1.Cells.Item(654,7).Value = "30/12/2025"

$ws1.Cells.Item(655,2).Value = "20:58:10"
$ws1.Cells.Item(655,3).Value = "21:32"
$ws1.Cells.Item(655,4).Value = "23_HERNANDEZ"
$ws1.Cells.Item(655,5).Value = 34
$ws1.Cells.Item(655,6).Value = "LP1912"
$ws1.Cells.Item(655,7).Value = "30/12/2025"

$ws1.Cells.Item(656,2).Value = "20:58:10"
$ws1.Cells.Item(656,3).Value = "21:38"
$ws1.Cells.Item(656,4).Value = "17_ROMERO"
$ws1.Cells.Item(656,5).Value = 40
$ws1.Cells.Item(656,6).Value = "LP1912"
$ws1.Cells.Item(656,7).Value = "30/12/2025"

$ws1.Cells.Item(657,2).Value = "20:58:10"
$ws1.Cells.Item(657,3).Value = "21:47"
$ws1.Cells.Item(657,4).Value = "16_SANTA ANA"
$ws1.Cells.Item(657,5).Value = 49
$ws1.Cells.Item(657,6).Value = "LP1912"
$ws1.Cells.Item(657,7).Value = "30/12/2025"

$ws1.Cells.Item(658,2).Value = "20:58:10"
$ws1.Cells.Item(658,3).Value = "21:47"
$ws1.Cells.Item(658,4).Value = "215A_EL PATO"
$ws1.Cells.Item(658,5).Value = 49
$ws1.Cells.Item(658,6).Value = "LP1912"
$ws1.Cells.Item(658,7).Value = "30/12/2025"

$ws1.Cells.Item(659,2).Value = "20:58:10"
$ws1.Cells.Item(659,3).Value = "21:51"
$ws1.Cells.Item(659,4).Value = "10_OLMOS"
$ws1.Cells.Item(659,5).Value = 53
$ws1.Cells.Item(659,6).Value = "LP1912"
$ws1.Cells.Item(659,7).Value = "30/12/2025"

$ws1.Cells.Item(660,2).Value = "20:58:10"
$ws1.Cells.Item(660,3).Value = "22:08"
$ws1.Cells.Item(660,4).Value = "17_ROMERO"
$ws1.Cells.Item(660,5).Value = 70
$ws1.Cells.Item(660,6).Value = "LP1912"
$ws1.Cells.Item(660,7).Value = "30/12/2025"

$ws1.Cells.Item(661,2).Value = "20:58:10"
$ws1.Cells.Item(661,3).Value = "22:23"
$ws1.Cells.Item(661,4).Value = "26_HERNANDEZ"
$ws1.Cells.Item(661,5).Value = 85
$ws1.Cells.Item(661,6).Value = "LP1912"
$ws1.Cells.Item(661,7).Value = "30/12/2025"

$ws1.Cells.Item(662,2).Value = "20:58:10"
$ws1.Cells.Item(662,3).Value = "22:25"
$ws1.Cells.Item(662,4).Value = "10_OLMOS"
$ws1.Cells.Item(662,5).Value = 87
$ws1.Cells.Item(662,6).Value = "LP1912"
$ws1.Cells.Item(662,7).Value = "30/12/2025"

$ws1.Cells.Item(663,2).Value = "20:58:10"
$ws1.Cells.Item(663,3).Value = "22:28"
$ws1.Cells.Item(663,4).Value = "84_COLONIA URQUIZA-ESC 49"
$ws1.Cells.Item(663,5).Value = 90
$ws1.Cells.Item(663,6).Value = "LP1912"
$ws1.Cells.Item(663,7).Value = "30/12/2025"

# ---- Sheet 2: LP1912-215 ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 20:58:21"
$ws2.Cells.Item(3,1).Value = "Total filas: 50"

$ws2.Cells.Item(50,2).Value = "30/12/2025"
$ws2.Cells.Item(50,3).Value = "20:58:10"
$ws2.Cells.Item(50,4).Value = "21:07"
$ws2.Cells.Item(50,5).Value = "215B_EL PATO"
$ws2.Cells.Item(50,6).Value = 9
$ws2.Cells.Item(50,7).Value = "LP1912"

$ws2.Cells.Item(51,2).Value = "30/12/2025"
$ws2.Cells.Item(51,3).Value = "20:58:10"
$ws2.Cells.Item(51,4).Value = "21:47"
$ws2.Cells.Item(51,5).Value = "215A_EL PATO"
$ws2.Cells.Item(51,6).Value = 49
$ws2.Cells.Item(51,7).Value = "LP1912"

# ---- Sheet 3: 6203-6173 ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Cells.Item(2,1).Value = "Última actualización: 30/12/2025 20:58:21"
$ws3.Cells.Item(3,1).Value = "Total filas: 82"

$ws3.Cells.Item(81,2).Value = "30/12/2025"
$ws3.Cells.Item(81,3).Value = "20:58:15"
$ws3.Cells.Item(81,4).Value = "21:29"
$ws3.Cells.Item(81,5).Value = "215C_LA PLATA"
$ws3.Cells.Item(81,6).Value = 31
$ws3.Cells.Item(81,7).Value = "L6203"

$ws3.Cells.Item(82,2).Value = "30/12/2025"
$ws3.Cells.Item(82,3).Value = "20:58:20"
$ws3.Cells.Item(82,4).Value = "22:05"
$ws3.Cells.Item(82,5).Value = "215A_LA PLATA"
$ws3.Cells.Item(82,6).Value = 67
$ws3.Cells.Item(82,7).Value = "L6173"

$ws3.Cells.Item(83,2).Value = "30/12/2025"
$ws3.Cells.Item(83,3).Value = "20:58:20"
$ws3.Cells.Item(83,4).Value = "22:21"
$ws3.Cells.Item(83,5).Value = "215B_LP-P MOR-40 Y 115"
$ws3.Cells.Item(83,6).Value = 83
$ws3.Cells.Item(83,7).Value = "L6173"
